# Update ARROIO_GRANDE.xlsx: rename two tabs and remove the
# "Desarquivamentos Pendentes" sheet that is no longer needed.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Drop the "Desarquivamentos Pendentes" tab entirely.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
